$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.534.22'
$ws.Range("E2").Value = '  +0.14%  '

# Row 3
$ws.Range("D3").Value = '1.846.94'
$ws.Range("E3").Value = '  -0.24%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '264.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.17%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.04%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5232'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.93%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3226'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.34%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06812'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.52%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.82'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.01%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7791'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.91%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07768'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.32%  '

# Row 13
$ws.Range("D13").Value = '1.858.64'
$ws.Range("E13").Value = '  +0.17%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.33%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.025'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.39%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9997'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.04%  '

# Row 17
$ws.Range("E17").Value = '  -1.29%  '

# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007981'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.87%  '

# Row 19
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.03%  '

# Row 20
$ws.Range("D20").Value = '26.553.62'
$ws.Range("E20").Value = '  +0.04%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.634'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.98%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.475'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.78%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.011'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.23%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '142.80'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.21%  '

# Row 25
$ws.Range("E25").Value = '  -7.29%  '

# Row 26
$ws.Range("E26").Value = '  +1.84%  '

# Row 27
$ws.Range("E27").Value = '  +0.02%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '111.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.02%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.178'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.98%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.122'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.49%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08730'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.30%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04840'
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.131'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.38%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7185'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.22%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.110'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.36%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01789'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.35%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.235'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.37%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.4874'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.07%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8993'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.66%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '110.62'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.93%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.027'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.53%  '

# Row 43
$ws.Range("E43").Value = '  -0.03%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.642'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.20%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4208'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.49%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05887'
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.046'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.67%  '

# Row 48
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.21%  '

# Row 49
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1224'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.29%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.8895'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.52%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.96'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.96%  '
